# "hero skill icon fix"
# The Icon column (F) of the HeroSkill sheet held ad-hoc pinyin names
# (chengjie, penhuo, dadun, zhiliao, miaozhun, shibei, yexing2, cisha).
# Replace them with the new standard skillN icon keys, row by row, while
# leaving every other column (Name/Des/etc.) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = "skill1"
$ws.Range("F5").Value  = "skill2"
$ws.Range("F6").Value  = "skill3"
$ws.Range("F7").Value  = "skill4"
$ws.Range("F8").Value  = "skill5"
$ws.Range("F9").Value  = "skill6"
$ws.Range("F10").Value = "skill7"
$ws.Range("F11").Value = "skill8"

# Match the author's final cursor position recorded in the sheet view.
$ws.Range("F7").Select()

# The workbook's custom "window background" theme tint (a pale green,
# CAEACD) reverts back to the Office-default white on this resave.
$themeColors = $wb.Theme.ThemeColorScheme
$themeColors.Item(2).RGB = 16777215
